# Auto-generated Excel COM-interop script to apply market-price / profit updates
# across the Gilgamesh Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2110.5334
$ws.Range("J17").Value = 2061.3572
$ws.Range("L17").Value = 6184.071599999999
$ws.Range("N17").Value = -6520.071599999999

$ws.Range("H41").Value = 1468.3846
$ws.Range("I41").Value = 309.6
$ws.Range("K41").Value = 309.6
$ws.Range("M41").Value = 130.4

$ws.Range("H46").Value = 2996
$ws.Range("J46").Value = 2996
$ws.Range("L46").Value = 8988
$ws.Range("N46").Value = -9226

$ws.Range("H60").Value = 2996
$ws.Range("J60").Value = 2996
$ws.Range("L60").Value = 8988
$ws.Range("N60").Value = -9956

$ws.Range("H64").Value = 250013120
$ws.Range("I64").Value = 17496.334
$ws.Range("K64").Value = 17496.334
$ws.Range("M64").Value = -17248.334

$ws.Range("H67").Value = 250013120
$ws.Range("I67").Value = 17496.334
$ws.Range("K67").Value = 17496.334
$ws.Range("M67").Value = -16638.334

$ws.Range("H69").Value = 7833
$ws.Range("J69").Value = 7833
$ws.Range("L69").Value = 23499
$ws.Range("N69").Value = -25247

$ws.Range("H72").Value = 7833
$ws.Range("J72").Value = 7833
$ws.Range("L72").Value = 70497
$ws.Range("N72").Value = -79233

$ws.Range("H74").Value = 18350.143
$ws.Range("I74").Value = 18350.143
$ws.Range("K74").Value = 18350.143
$ws.Range("M74").Value = -17414.143

$ws.Range("H75").Value = 125000
$ws.Range("J75").Value = 125000
$ws.Range("L75").Value = 125000
$ws.Range("N75").Value = -126872

$ws.Range("H77").Value = 18350.143
$ws.Range("I77").Value = 18350.143
$ws.Range("K77").Value = 91750.715
$ws.Range("M77").Value = -87070.715

$ws.Range("H78").Value = 125000
$ws.Range("J78").Value = 125000
$ws.Range("L78").Value = 375000
$ws.Range("N78").Value = -384360

$ws.Range("H113").Value = 1300
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -7708

$ws.Range("H137").Value = 11234.25
$ws.Range("J137").Value = 17143.428
$ws.Range("L137").Value = 51430.284
$ws.Range("N137").Value = -56530.284

$ws.Range("H138").Value = 325129.28
$ws.Range("J138").Value = 377401.8
$ws.Range("L138").Value = 1132205.4
$ws.Range("N138").Value = -1142485.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1586
$ws.Range("I2").Value = 1591.8
$ws.Range("K2").Value = 1591.8
$ws.Range("M2").Value = -1478.8

$ws.Range("H32").Value = 3872.9126
$ws.Range("I32").Value = 3491.3376
$ws.Range("K32").Value = 3491.3376
$ws.Range("M32").Value = -3204.3376

$ws.Range("H63").Value = 6334.8335
$ws.Range("I63").Value = 4752.25
$ws.Range("K63").Value = 4752.25
$ws.Range("M63").Value = -4066.25

$ws.Range("H66").Value = 6334.8335
$ws.Range("I66").Value = 4752.25
$ws.Range("K66").Value = 23761.25
$ws.Range("M66").Value = -20329.25

$ws.Range("H116").Value = 1586
$ws.Range("I116").Value = 1591.8
$ws.Range("K116").Value = 1591.8
$ws.Range("M116").Value = 702.2

$ws.Range("H132").Value = 2098.2
$ws.Range("I132").Value = 1755.2667
$ws.Range("K132").Value = 5265.800099999999
$ws.Range("M132").Value = -2735.800099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1586
$ws.Range("I3").Value = 1591.8
$ws.Range("K3").Value = 1591.8
$ws.Range("M3").Value = -1477.8

$ws.Range("H81").Value = 38491.5
$ws.Range("J81").Value = 38491.5
$ws.Range("L81").Value = 38491.5
$ws.Range("N81").Value = -40613.5

$ws.Range("H82").Value = 81072.57000000001
$ws.Range("I82").Value = 72910.75
$ws.Range("J82").Value = 91955
$ws.Range("K82").Value = 72910.75
$ws.Range("L82").Value = 91955
$ws.Range("M82").Value = -72527.75
$ws.Range("N82").Value = -92721

$ws.Range("H84").Value = 38491.5
$ws.Range("J84").Value = 38491.5
$ws.Range("L84").Value = 115474.5
$ws.Range("N84").Value = -126082.5

$ws.Range("H85").Value = 81072.57000000001
$ws.Range("I85").Value = 72910.75
$ws.Range("J85").Value = 91955
$ws.Range("K85").Value = 72910.75
$ws.Range("L85").Value = 91955
$ws.Range("M85").Value = -71584.75
$ws.Range("N85").Value = -94607

$ws.Range("H97").Value = 11122.25
$ws.Range("I97").Value = 11122.25
$ws.Range("K97").Value = 11122.25
$ws.Range("M97").Value = -10131.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4000.8
$ws.Range("I31").Value = 3668.2903
$ws.Range("J31").Value = 4303.9707
$ws.Range("K31").Value = 3668.2903
$ws.Range("L31").Value = 4303.9707
$ws.Range("M31").Value = -3373.2903
$ws.Range("N31").Value = -4893.9707

$ws.Range("H34").Value = 4000.8
$ws.Range("I34").Value = 3668.2903
$ws.Range("J34").Value = 4303.9707
$ws.Range("K34").Value = 3668.2903
$ws.Range("L34").Value = 4303.9707
$ws.Range("M34").Value = -3466.2903
$ws.Range("N34").Value = -4707.9707

$ws.Range("H62").Value = 7149695
$ws.Range("I62").Value = 20004140
$ws.Range("K62").Value = 20004140
$ws.Range("M62").Value = -20003516

$ws.Range("H65").Value = 7149695
$ws.Range("I65").Value = 20004140
$ws.Range("K65").Value = 100020700
$ws.Range("M65").Value = -100017580

$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50629

$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -52183

$ws.Range("H122").Value = 2031.52
$ws.Range("I122").Value = 2049.3333
$ws.Range("J122").Value = 1985.7142
$ws.Range("K122").Value = 6147.999899999999
$ws.Range("L122").Value = 5957.142599999999
$ws.Range("M122").Value = -3697.999899999999
$ws.Range("N122").Value = -10857.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 397
$ws.Range("I98").Value = 641
$ws.Range("J98").Value = 234.33333
$ws.Range("K98").Value = 1923
$ws.Range("L98").Value = 702.99999
$ws.Range("M98").Value = -425
$ws.Range("N98").Value = -3698.99999

$ws.Range("H115").Value = 305540
$ws.Range("I115").Value = 2500
$ws.Range("J115").Value = 339211.12
$ws.Range("K115").Value = 7500
$ws.Range("L115").Value = 1017633.36
$ws.Range("M115").Value = -6325
$ws.Range("N115").Value = -1019983.36

$ws.Range("H125").Value = 6000
$ws.Range("J125").Value = 6000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -27840

$ws.Range("H129").Value = 2155.8572
$ws.Range("J129").Value = 1815.5
$ws.Range("L129").Value = 5446.5
$ws.Range("N129").Value = -15446.5

$ws.Range("H134").Value = 4837.7
$ws.Range("I134").Value = 2696.8572
$ws.Range("J134").Value = 9833
$ws.Range("K134").Value = 8090.571599999999
$ws.Range("L134").Value = 29499
$ws.Range("M134").Value = -3020.571599999999
$ws.Range("N134").Value = -39639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142861570
$ws.Range("J80").Value = 6333.3335
$ws.Range("L80").Value = 6333.3335
$ws.Range("N80").Value = -8329.333500000001

$ws.Range("H83").Value = 142861570
$ws.Range("J83").Value = 6333.3335
$ws.Range("L83").Value = 31666.6675
$ws.Range("N83").Value = -41650.6675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5249.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 5249.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H136").Value = 9971.444
$ws.Range("I136").Value = 13497
$ws.Range("K136").Value = 40491
$ws.Range("M136").Value = -37941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47621
$ws.Range("J70").Value = 47621
$ws.Range("L70").Value = 47621
$ws.Range("N70").Value = -48251

$ws.Range("H73").Value = 47621
$ws.Range("J73").Value = 47621
$ws.Range("L73").Value = 47621
$ws.Range("N73").Value = -49805
